# refactor : example.xlsx file
#
# The dates that used to live in row 5 (A5:E5) as real date serial values
# (formatted with a date number format) are converted to plain text
# values holding the ISO date strings "2023-03-20" .. "2023-03-24".
# The columns are widened to fit the new text, and the active selection
# moves to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: turn the date values into literal text dates ------------------
# NumberFormat has to be switched to Text ("@", Excel's built-in format 49)
# BEFORE the string values are written, otherwise Excel would keep storing
# them as date serial numbers.
$dateRange = $ws.Range("A5:E5")
$dateRange.NumberFormat = "@"

$ws.Range("A5").Value = "2023-03-20"
$ws.Range("B5").Value = "2023-03-21"
$ws.Range("C5").Value = "2023-03-22"
$ws.Range("D5").Value = "2023-03-23"
$ws.Range("E5").Value = "2023-03-24"

# --- Widen columns A:E so the new text values fit nicely -------------------
$ws.Range("A1:E5").EntireColumn.ColumnWidth = 10.4

# --- Move the active selection to C6 ---------------------------------------
$ws.Range("C6").Select() | Out-Null
